# EL-NumberValidation.xlsx / "NumberError" sheet
#
# This replays a fresh Katalon RAD test-suite run over the 46 "Existing
# Liability w/Notice Number" number-validation test cases (rows 2-47):
#   - Column A ("Result") is (re)stamped "Pass" for every test case, including
#     16 rows (32-47) that previously only had the static Tax Type / expected
#     error-message columns (C-G) populated and had never recorded a
#     Result/Date before.
#   - Column B ("Date") is (re)stamped with the timestamp at which each test
#     case executed during this run (Sun Jan 14 2024, ~16:16-16:23 EST).
#
# Columns C-G (TaxType/TaxTypeEL/NotiInvoNumber/ErrorMessage) already hold the
# correct static test-case data and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$resultDates = @(
    "Sun Jan 14 16:16:29 EST 2024",
    "Sun Jan 14 16:16:38 EST 2024",
    "Sun Jan 14 16:16:47 EST 2024",
    "Sun Jan 14 16:16:56 EST 2024",
    "Sun Jan 14 16:17:06 EST 2024",
    "Sun Jan 14 16:17:16 EST 2024",
    "Sun Jan 14 16:17:25 EST 2024",
    "Sun Jan 14 16:17:34 EST 2024",
    "Sun Jan 14 16:17:43 EST 2024",
    "Sun Jan 14 16:17:52 EST 2024",
    "Sun Jan 14 16:18:01 EST 2024",
    "Sun Jan 14 16:18:11 EST 2024",
    "Sun Jan 14 16:18:20 EST 2024",
    "Sun Jan 14 16:18:29 EST 2024",
    "Sun Jan 14 16:18:38 EST 2024",
    "Sun Jan 14 16:18:48 EST 2024",
    "Sun Jan 14 16:18:57 EST 2024",
    "Sun Jan 14 16:19:06 EST 2024",
    "Sun Jan 14 16:19:16 EST 2024",
    "Sun Jan 14 16:19:25 EST 2024",
    "Sun Jan 14 16:19:34 EST 2024",
    "Sun Jan 14 16:19:43 EST 2024",
    "Sun Jan 14 16:19:52 EST 2024",
    "Sun Jan 14 16:20:01 EST 2024",
    "Sun Jan 14 16:20:10 EST 2024",
    "Sun Jan 14 16:20:20 EST 2024",
    "Sun Jan 14 16:20:29 EST 2024",
    "Sun Jan 14 16:20:38 EST 2024",
    "Sun Jan 14 16:20:47 EST 2024",
    "Sun Jan 14 16:20:56 EST 2024",
    "Sun Jan 14 16:21:06 EST 2024",
    "Sun Jan 14 16:21:15 EST 2024",
    "Sun Jan 14 16:21:24 EST 2024",
    "Sun Jan 14 16:21:33 EST 2024",
    "Sun Jan 14 16:21:42 EST 2024",
    "Sun Jan 14 16:21:51 EST 2024",
    "Sun Jan 14 16:22:00 EST 2024",
    "Sun Jan 14 16:22:11 EST 2024",
    "Sun Jan 14 16:22:20 EST 2024",
    "Sun Jan 14 16:22:29 EST 2024",
    "Sun Jan 14 16:22:38 EST 2024",
    "Sun Jan 14 16:22:47 EST 2024",
    "Sun Jan 14 16:22:56 EST 2024",
    "Sun Jan 14 16:23:05 EST 2024",
    "Sun Jan 14 16:23:15 EST 2024",
    "Sun Jan 14 16:23:24 EST 2024"
)

for ($i = 0; $i -lt $resultDates.Length; $i++) {
    $row = $i + 2

    $resultCell = $ws.Cells.Item($row, 1)
    $resultCell.Value = "Pass"
    # Rows 32-47 never had a Result/Date before; make sure the brand-new
    # cells pick up the sheet's default (unstyled) look instead of
    # inheriting the column's border/format style, matching the rest of
    # the Result column.
    $resultCell.Style = "Normal"

    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.Value = $resultDates[$i]
    $dateCell.Style = "Normal"
}
